# 19/10/2023
# * UT_5023x updated (closed loop simulations)
# * UT_5030x added (forward crosstalk simulations)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New test rows (UT_5030a..UT_5030e) appended after the existing data (row 33).
$newRows = @(
    @{ Row = 34; A = "DRE_DMX_UT_5030a"; B = "Simulation of forward crosstalk with pulse shaping=1" + [char]10 + "and SAMPLING_DELAY=0x0008"; C = "PASS"; D = "Crosstalk is observed" },
    @{ Row = 35; A = "DRE_DMX_UT_5030b"; B = "Simulation of forward crosstalk with pulse shaping=1" + [char]10 + "and SAMPLING_DELAY=x000B";  C = "PASS"; D = "Crosstalk is observed" },
    @{ Row = 36; A = "DRE_DMX_UT_5030c"; B = "Simulation of forward crosstalk with pulse shaping=1" + [char]10 + "and SAMPLING_DELAY=0x0012"; C = "PASS"; D = "No crosstalk " },
    @{ Row = 37; A = "DRE_DMX_UT_5030d"; B = "Simulation of forward crosstalk with pulse shaping=3" + [char]10 + "and SAMPLING_DELAY=0x0008"; C = "PASS"; D = "Crosstalk is observed" },
    @{ Row = 38; A = "DRE_DMX_UT_5030e"; B = "Simulation of forward crosstalk with pulse shaping=3" + [char]10 + "and SAMPLING_DELAY=0x000B"; C = "PASS"; D = "Crosstalk is observed" }
)

# Column A and D for the first few rows, then the (wrapped) column B text, then
# the remaining cells - following the same cell fill order used when the rows
# were originally entered.
$ws.Cells.Item(34, 1).Value = $newRows[0].A
$ws.Cells.Item(34, 4).Value = $newRows[0].D
$ws.Cells.Item(35, 1).Value = $newRows[1].A
$ws.Cells.Item(36, 1).Value = $newRows[2].A
$ws.Cells.Item(34, 2).Value = $newRows[0].B
$ws.Cells.Item(35, 2).Value = $newRows[1].B
$ws.Cells.Item(36, 2).Value = $newRows[2].B
$ws.Cells.Item(37, 1).Value = $newRows[3].A
$ws.Cells.Item(37, 2).Value = $newRows[3].B
$ws.Cells.Item(38, 1).Value = $newRows[4].A
$ws.Cells.Item(38, 2).Value = $newRows[4].B
$ws.Cells.Item(36, 4).Value = $newRows[2].D

foreach ($r in $newRows) {
    $rowIdx = $r.Row

    $ws.Cells.Item($rowIdx, 3).Value = $r.C
    $ws.Cells.Item($rowIdx, 4).Value = $r.D

    $bCell = $ws.Cells.Item($rowIdx, 2)
    $bCell.WrapText = $true

    $ws.Rows.Item($rowIdx).RowHeight = 30
}

# Scroll the view down and select the final block, mirroring the author's
# on-screen state when the edit was made.
$excel.ActiveWindow.ScrollRow = 22
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D37:D38").Select()
